$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet "Binary Search 2": scroll the frozen view so the visible
# top row moves from 8 to 7, and move the selection from F11 to F7 ---
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws2.Range("F7").Select()

# --- Sheet "Notes": add a new row (5) describing the trailing-zeroes
# note, then leave it as the active / selected sheet ---
$ws3.Range("B6").Value = 5

$ws3.Range("D6").Value = "No. of trailing 0s in N! [N factorial]"
$ws3.Range("D6").WrapText = $true

$ws3.Hyperlinks.Add($ws3.Range("E6"), "https://takeuforward.org/data-structure/count-trailing-zeroes/", "", "", "Count Trailing Zeroes - Tutorial [Updated] (takeuforward.org)")
$ws3.Range("E6").Style = "Hyperlink"
$ws3.Range("E6").WrapText = $true

$ws3.Range("C6").Value = "Notes 3"

$ws3.Range("G6").Value = "Good problem to understand logrithmic complexity"
$ws3.Range("G6").WrapText = $true

# New column G needs an explicit width of 18 (characters)
$ws3.Columns.Item(7).ColumnWidth = 17.166666666666668

# Make "Notes" the active sheet/tab, and select the newly added data cell
$ws3.Activate()
$ws3.Range("H6").Select()
